$wb = $excel.ActiveWorkbook

# --- ALC row 40 (hunk @@ -2622,25 +2622,25 @@) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3408.8
$ws.Range("J40").Value = 3509.7778
$ws.Range("L40").Value = 3509.7778
$ws.Range("N40").Value = -3859.7778

# --- ALC row 99 (hunk @@ -5591,25 +5591,25 @@) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 1275.8182
$ws.Range("I99").Value = 219.14285
$ws.Range("J99").Value = 3125
$ws.Range("K99").Value = 657.4285500000001
$ws.Range("L99").Value = 9375
$ws.Range("M99").Value = 840.5714499999999
$ws.Range("N99").Value = -12371

# --- ALC row 113 (hunk @@ -6295,25 +6295,25 @@) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3452.476
$ws.Range("I113").Value = 2362.5
$ws.Range("J113").Value = 3708.9412
$ws.Range("K113").Value = 2362.5
$ws.Range("L113").Value = 3708.9412
$ws.Range("M113").Value = 891.5
$ws.Range("N113").Value = -10216.9412

# --- ALC row 115 (hunk @@ -6396,25 +6396,25 @@) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 3657.8
$ws.Range("I115").Value = 1328.3334
$ws.Range("J115").Value = 7152
$ws.Range("K115").Value = 3985.0002
$ws.Range("L115").Value = 21456
$ws.Range("M115").Value = -2418.0002
$ws.Range("N115").Value = -24590

# --- ALC row 127 (hunk @@ -6987,25 +6987,25 @@) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 780.38464
$ws.Range("I127").Value = 638.3333
$ws.Range("J127").Value = 1100
$ws.Range("K127").Value = 1914.9999
$ws.Range("L127").Value = 3300
$ws.Range("M127").Value = 3045.0001
$ws.Range("N127").Value = -13220

# --- ALC row 132 (hunk @@ -7244,25 +7244,25 @@) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6153.5317
$ws.Range("I132").Value = 4497.846
$ws.Range("J132").Value = 14225
$ws.Range("K132").Value = 13493.538
$ws.Range("L132").Value = 42675
$ws.Range("M132").Value = -10963.538
$ws.Range("N132").Value = -47735

# --- ALC row 137 (hunk @@ -7492,25 +7492,25 @@) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2090
$ws.Range("I137").Value = 2452.7778
$ws.Range("J137").Value = 1364.4445
$ws.Range("K137").Value = 7358.3334
$ws.Range("L137").Value = 4093.3335
$ws.Range("M137").Value = -4808.3334
$ws.Range("N137").Value = -9193.333500000001

# --- ARM row 61 (hunk @@ -10716,25 +10716,25 @@) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3482.0625
$ws.Range("I61").Value = 1958
$ws.Range("J61").Value = 4667.4443
$ws.Range("K61").Value = 1958
$ws.Range("L61").Value = 4667.4443
$ws.Range("M61").Value = -1746
$ws.Range("N61").Value = -5091.4443

# --- ARM row 74 (hunk @@ -11350,22 +11350,22 @@) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1420.8695
$ws.Range("I74").Value = 1287.0286
$ws.Range("K74").Value = 1287.0286
$ws.Range("M74").Value = -413.0286000000001

# --- ARM row 77 (hunk @@ -11497,22 +11497,22 @@) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1420.8695
$ws.Range("I77").Value = 1287.0286
$ws.Range("K77").Value = 6435.143
$ws.Range("M77").Value = -2067.143

# --- ARM row 122 (hunk @@ -13690,25 +13690,25 @@) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2376.2068
$ws.Range("I122").Value = 1226.2307
$ws.Range("J122").Value = 12342.667
$ws.Range("K122").Value = 3678.6921
$ws.Range("L122").Value = 37028.001
$ws.Range("M122").Value = -1228.6921
$ws.Range("N122").Value = -41928.001

# --- ARM row 136 (hunk @@ -14370,25 +14370,25 @@) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3482.0625
$ws.Range("I136").Value = 1958
$ws.Range("J136").Value = 4667.4443
$ws.Range("K136").Value = 5874
$ws.Range("L136").Value = 14002.3329
$ws.Range("M136").Value = -3324
$ws.Range("N136").Value = -19102.3329

# --- BSM row 134 (hunk @@ -21184,25 +21184,25 @@) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8571.174000000001
$ws.Range("I134").Value = 3248.6
$ws.Range("J134").Value = 18551
$ws.Range("K134").Value = 9745.799999999999
$ws.Range("L134").Value = 55653
$ws.Range("M134").Value = -7210.799999999999
$ws.Range("N134").Value = -60723

# --- CRP row 6 (hunk @@ -21875,25 +21875,25 @@) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 22000840
$ws.Range("I6").Value = 55000000
$ws.Range("J6").Value = 1400
$ws.Range("K6").Value = 55000000
$ws.Range("L6").Value = 1400
$ws.Range("M6").Value = -54999887
$ws.Range("N6").Value = -1626

# --- CRP row 7 (hunk @@ -21927,25 +21927,25 @@) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 4645.7393
$ws.Range("I7").Value = 9116.091
$ws.Range("J7").Value = 547.9167
$ws.Range("K7").Value = 9116.091
$ws.Range("L7").Value = 547.9167
$ws.Range("M7").Value = -9003.091
$ws.Range("N7").Value = -773.9167

# --- CRP row 17 (hunk @@ -22423,25 +22423,22 @@) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 4800
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 4800
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4800
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -5148

# --- CRP row 25 (hunk @@ -22833,25 +22830,25 @@) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 31392.125
$ws.Range("I25").Value = 11111
$ws.Range("J25").Value = 34289.43
$ws.Range("K25").Value = 11111
$ws.Range("L25").Value = 34289.43
$ws.Range("M25").Value = -10937
$ws.Range("N25").Value = -34637.43

# --- CRP row 31 (hunk @@ -23130,22 +23127,25 @@) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9526061
$ws.Range("I31").Value = 1823.3448
$ws.Range("J31").Value = 55559876
$ws.Range("K31").Value = 1823.3448
$ws.Range("L31").Value = 55559876
$ws.Range("M31").Value = -1528.3448
$ws.Range("N31").Value = -55560466

# --- CRP row 34 (hunk @@ -23277,22 +23277,25 @@) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 9526061
$ws.Range("I34").Value = 1823.3448
$ws.Range("J34").Value = 55559876
$ws.Range("K34").Value = 1823.3448
$ws.Range("L34").Value = 55559876
$ws.Range("M34").Value = -1621.3448
$ws.Range("N34").Value = -55560280

# --- CRP row 140 (hunk @@ -28489,22 +28492,22 @@) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 44390
$ws.Range("J140").Value = 44390
$ws.Range("L140").Value = 44390
$ws.Range("N140").Value = -54750

# --- CUL row 5 (hunk @@ -28837,22 +28840,22 @@) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 609.05554
$ws.Range("I5").Value = 357.82352
$ws.Range("K5").Value = 1073.47056
$ws.Range("M5").Value = -961.47056

# --- CUL row 7 (hunk @@ -28941,25 +28944,25 @@) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 610
$ws.Range("I7").Value = 350
$ws.Range("J7").Value = 783.3333
$ws.Range("K7").Value = 1050
$ws.Range("L7").Value = 2349.9999
$ws.Range("M7").Value = -938
$ws.Range("N7").Value = -2573.9999

# --- CUL row 36 (hunk @@ -30425,22 +30428,22 @@) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 2425
$ws.Range("I36").Value = 566.6667
$ws.Range("K36").Value = 1700.0001
$ws.Range("M36").Value = -1531.0001

# --- CUL row 113 (hunk @@ -34345,25 +34348,25 @@) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 702.4737
$ws.Range("I113").Value = 661.6667
$ws.Range("J113").Value = 721.3077
$ws.Range("K113").Value = 1985.0001
$ws.Range("L113").Value = 2163.9231
$ws.Range("M113").Value = 184.9999
$ws.Range("N113").Value = -6503.9231

# --- CUL row 131 (hunk @@ -35275,25 +35278,25 @@) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 239285.72
$ws.Range("I131").Value = 1250412.1
$ws.Range("J131").Value = 1373.6177
$ws.Range("K131").Value = 3751236.3
$ws.Range("L131").Value = 4120.8531
$ws.Range("M131").Value = -3746196.3
$ws.Range("N131").Value = -14200.8531

# --- CUL row 135 (hunk @@ -35483,22 +35486,22 @@) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 609.05554
$ws.Range("I135").Value = 357.82352
$ws.Range("K135").Value = 3220.41168
$ws.Range("M135").Value = -685.4116799999997

# --- GSM row 122 (hunk @@ -41767,25 +41770,25 @@) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 224675.56
$ws.Range("I122").Value = 346897.6
$ws.Range("J122").Value = 3148.125
$ws.Range("K122").Value = 1040692.8
$ws.Range("L122").Value = 9444.375
$ws.Range("M122").Value = -1038242.8
$ws.Range("N122").Value = -14344.375

# --- LTW row 7 (hunk @@ -43089,22 +43092,25 @@) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2142.25
$ws.Range("I7").Value = 1982.5454
$ws.Range("J7").Value = 3899
$ws.Range("K7").Value = 1982.5454
$ws.Range("L7").Value = 3899
$ws.Range("M7").Value = -1870.5454
$ws.Range("N7").Value = -4123

# --- LTW row 61 (hunk @@ -45732,25 +45738,25 @@) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4662.727
$ws.Range("I61").Value = 1125
$ws.Range("J61").Value = 6684.2856
$ws.Range("K61").Value = 1125
$ws.Range("L61").Value = 6684.2856
$ws.Range("M61").Value = -923
$ws.Range("N61").Value = -7088.2856

# --- LTW row 113 (hunk @@ -48259,25 +48265,25 @@) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4662.727
$ws.Range("I113").Value = 1125
$ws.Range("J113").Value = 6684.2856
$ws.Range("K113").Value = 1125
$ws.Range("L113").Value = 6684.2856
$ws.Range("M113").Value = 1045
$ws.Range("N113").Value = -11024.2856

# --- LTW row 126 (hunk @@ -48893,22 +48899,25 @@) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2142.25
$ws.Range("I126").Value = 1982.5454
$ws.Range("J126").Value = 3899
$ws.Range("K126").Value = 5947.6362
$ws.Range("L126").Value = 3899
$ws.Range("M126").Value = -3477.6362
$ws.Range("N126").Value = -16637

# --- WVR row 100 (hunk @@ -54588,25 +54597,25 @@) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 614.1429000000001
$ws.Range("I100").Value = 522.05554
$ws.Range("J100").Value = 1166.6666
$ws.Range("K100").Value = 1044.11108
$ws.Range("L100").Value = 2333.3332
$ws.Range("M100").Value = -503.1110799999999
$ws.Range("N100").Value = -3415.3332

# --- WVR row 113 (hunk @@ -55222,25 +55231,25 @@) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I113").Value = 358.25
$ws.Range("J113").Value = 1079.3889
$ws.Range("K113").Value = 1074.75
$ws.Range("L113").Value = 3238.1667
$ws.Range("M113").Value = 1095.25
$ws.Range("N113").Value = -7578.1667
